$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "321.92"
Set-TextValue "E2" "8.45%"
Set-TextValue "D3" "51.13"
Set-TextValue "E3" "21.97%"
Set-TextValue "D4" "5.347"
Set-TextValue "E4" "6.85%"
Set-TextValue "D6" "4.577"
Set-TextValue "E6" "4.54%"
Set-TextValue "D7" "1.652"
Set-TextValue "E7" "4.77%"
Set-TextValue "D8" "1.119"
Set-TextValue "E8" "20.65%"
Set-TextValue "D9" "0.1317"
Set-TextValue "E9" "10.25%"
Set-TextValue "E10" "7.33%"
Set-TextValue "D11" "0.09646"
Set-TextValue "E11" "8.90%"
Set-TextValue "D12" "0.04574"
Set-TextValue "E12" "11.88%"
Set-TextValue "E13" "-0.14%"
Set-TextValue "D14" "0.001317"
Set-TextValue "E14" "2.25%"
Set-TextValue "D15" "0.005805"
Set-TextValue "E15" "-3.65%"
Set-TextValue "D16" "3.378"
Set-TextValue "D17" "2.432"
Set-TextValue "E17" "1.28%"
Set-TextValue "D18" "0.3395"
Set-TextValue "E18" "2.44%"
Set-TextValue "D19" "8.193"
Set-TextValue "E19" "1.45%"
Set-TextValue "E20" "0.49%"
Set-TextValue "E21" "-11.42%"
Set-TextValue "D22" "0.04307"
Set-TextValue "E22" "5.00%"
Set-TextValue "E23" "3.10%"
Set-TextValue "D24" "0.004308"
Set-TextValue "E24" "10.79%"
Set-TextValue "D25" "0.0001347"
Set-TextValue "E25" "9.38%"
Set-TextValue "E26" "-0.24%"
Set-TextValue "D38" "0.02769"
Set-TextValue "E38" "14.77%"
Set-TextValue "D39" "0.05541"
Set-TextValue "E39" "6.42%"
Set-TextValue "D40" "0.006288"
Set-TextValue "E40" "-0.33%"
Set-TextValue "D41" "0.007779"
Set-TextValue "E41" "-0.91%"
Set-TextValue "E42" "9.02%"
Set-TextValue "D43" "0.007676"
Set-TextValue "E43" "3.83%"
Set-TextValue "D44" "0.008813"
Set-TextValue "E44" "18.78%"
Set-TextValue "D45" "0.3529"
Set-TextValue "E45" "19.63%"
Set-TextValue "D46" "0.00006821"
Set-TextValue "E46" "5.71%"
Set-TextValue "E47" "-0.33%"
Set-TextValue "D48" "0.06040"
Set-TextValue "E48" "78.74%"
Set-TextValue "E50" "-0.33%"
Set-TextValue "D51" "0.0001996"
Set-TextValue "E51" "-0.33%"
